{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic expressions (e.g. \"21+46=\"). The edit replaces the text of\n// every cell with a new expression, cell-for-cell in row-major order\n// (the same physical cell/run/formatting is kept - only the visible\n// text changes). Build the new grid of values (row-major) exactly as\n// it must appear after the edit, then push it back into the table in\n// one shot via the Table.values property, which rewrites each cell's\n// text while preserving the existing run/paragraph formatting.\n\nconst newValues = [\n  [\"86-45=\", \"54+27=\", \"28-18=\", \"72-31=\", \"22+28=\"],\n  [\"86-82=\", \"34+15=\", \"17+27=\", \"10-1=\", \"8+32=\"],\n  [\"28+48=\", \"89-80=\", \"65+25=\", \"16+83=\", \"15+80=\"],\n  [\"78-56=\", \"89-39=\", \"34-25=\", \"53-45=\", \"22+76=\"],\n  [\"34-11=\", \"4+43=\", \"81-24=\", \"29+53=\", \"86+11=\"],\n  [\"91-7=\", \"19-4=\", \"65-32=\", \"87-13=\", \"62-42=\"],\n  [\"19-10=\", \"39+13=\", \"51-29=\", \"36-12=\", \"60+8=\"],\n  [\"50-29=\", \"75-25=\", \"16+82=\", \"95-24=\", \"87-56=\"],\n  [\"94-69=\", \"95+2=\", \"22+4=\", \"62+35=\", \"94-75=\"],\n  [\"41-18=\", \"89-15=\", \"80-74=\", \"93-81=\", \"0+25=\"],\n  [\"15+23=\", \"34+50=\", \"19+38=\", \"80+7=\", \"9+60=\"],\n  [\"22+60=\", \"47+15=\", \"24+17=\", \"51+2=\", \"62-40=\"],\n  [\"39+35=\", \"21+21=\", \"50+38=\", \"28-18=\", \"94-4=\"],\n  [\"85+5=\", \"67+20=\", \"23+74=\", \"82-52=\", \"21+75=\"],\n  [\"26-3=\", \"97-20=\", \"41-33=\", \"4+48=\", \"0+33=\"],\n  [\"33+26=\", \"45+36=\", \"17-1=\", \"3+19=\", \"23+58=\"],\n  [\"77-55=\", \"87-62=\", \"94-69=\", \"22-6=\", \"48+30=\"],\n  [\"90+5=\", \"79+9=\", \"42+13=\", \"96-33=\", \"30+63=\"],\n  [\"21+38=\", \"4+87=\", \"63+32=\", \"31-0=\", \"35+7=\"],\n  [\"58-22=\", \"81-22=\", \"32+62=\", \"46-39=\", \"90-67=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Writing the whole grid back in one shot replaces each cell's text in\n// place (row-major order, matching the physical table layout) while\n// preserving the existing run/paragraph formatting of every cell - this\n// is safe even though a handful of the \"before\" expressions repeat at\n// multiple positions, because the values are applied positionally.\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# The document body holds a single 20-row x 5-column table of simple\n# arithmetic prompts (e.g. \"21+46=\"). The edit swaps the text of every\n# cell, in row-major reading order, for a new expression - same cell,\n# same formatting, only the visible text changes. $newValues below is\n# that row-major list of replacement strings (index 0 = row1/col1,\n# index 1 = row1/col2, ... index 99 = row20/col5).\n\n$newValues = @(\n    \"86-45=\",\n    \"54+27=\",\n    \"28-18=\",\n    \"72-31=\",\n    \"22+28=\",\n    \"86-82=\",\n    \"34+15=\",\n    \"17+27=\",\n    \"10-1=\",\n    \"8+32=\",\n    \"28+48=\",\n    \"89-80=\",\n    \"65+25=\",\n    \"16+83=\",\n    \"15+80=\",\n    \"78-56=\",\n    \"89-39=\",\n    \"34-25=\",\n    \"53-45=\",\n    \"22+76=\",\n    \"34-11=\",\n    \"4+43=\",\n    \"81-24=\",\n    \"29+53=\",\n    \"86+11=\",\n    \"91-7=\",\n    \"19-4=\",\n    \"65-32=\",\n    \"87-13=\",\n    \"62-42=\",\n    \"19-10=\",\n    \"39+13=\",\n    \"51-29=\",\n    \"36-12=\",\n    \"60+8=\",\n    \"50-29=\",\n    \"75-25=\",\n    \"16+82=\",\n    \"95-24=\",\n    \"87-56=\",\n    \"94-69=\",\n    \"95+2=\",\n    \"22+4=\",\n    \"62+35=\",\n    \"94-75=\",\n    \"41-18=\",\n    \"89-15=\",\n    \"80-74=\",\n    \"93-81=\",\n    \"0+25=\",\n    \"15+23=\",\n    \"34+50=\",\n    \"19+38=\",\n    \"80+7=\",\n    \"9+60=\",\n    \"22+60=\",\n    \"47+15=\",\n    \"24+17=\",\n    \"51+2=\",\n    \"62-40=\",\n    \"39+35=\",\n    \"21+21=\",\n    \"50+38=\",\n    \"28-18=\",\n    \"94-4=\",\n    \"85+5=\",\n    \"67+20=\",\n    \"23+74=\",\n    \"82-52=\",\n    \"21+75=\",\n    \"26-3=\",\n    \"97-20=\",\n    \"41-33=\",\n    \"4+48=\",\n    \"0+33=\",\n    \"33+26=\",\n    \"45+36=\",\n    \"17-1=\",\n    \"3+19=\",\n    \"23+58=\",\n    \"77-55=\",\n    \"87-62=\",\n    \"94-69=\",\n    \"22-6=\",\n    \"48+30=\",\n    \"90+5=\",\n    \"79+9=\",\n    \"42+13=\",\n    \"96-33=\",\n    \"30+63=\",\n    \"21+38=\",\n    \"4+87=\",\n    \"63+32=\",\n    \"31-0=\",\n    \"35+7=\",\n    \"58-22=\",\n    \"81-22=\",\n    \"32+62=\",\n    \"46-39=\",\n    \"90-67=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        # Assigning Range.Text replaces just the visible run text and\n        # keeps the cell's end-of-cell marker and existing run/paragraph\n        # formatting (font, size, alignment) intact.\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
